# Refresh the crypto price/volume table (prices in column D, 1h volume
# change in column E) to the latest scraped snapshot. Two rows also swap
# rank position with their neighbour (Cosmos/Monero, BinanceUSD/FTXToken,
# InjectiveProtocol/TrustWalletToken), so Coin name + Link + Price + Volume
# are rewritten together for those rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Many price cells look numeric (e.g. "227.62", "1.00") and plain
# Range.Value assignment would let Excel auto-coerce them to numbers,
# dropping significant trailing zeros / changing the cell type. Force text
# via NumberFormat "@" before writing, then borrow the (unmodified, default)
# style from an untouched row in the same column so we don't leave a
# leftover explicit number-format style on the cell.
function Set-TextCell($addr, $value, $styleDonor) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = $ws.Range($styleDonor).Style
}

Set-TextCell "D2" '39.620.04' "D23"
Set-TextCell "E2" '  +2.12%  ' "E23"
Set-TextCell "D3" '2.160.68' "D23"
Set-TextCell "E3" '  +2.66%  ' "E23"
Set-TextCell "E4" '  +0.44%  ' "E23"
Set-TextCell "D5" '227.62' "D23"
Set-TextCell "E5" '  -0.09%  ' "E23"
Set-TextCell "D6" '0.629' "D23"
Set-TextCell "E6" '  +2.11%  ' "E23"
Set-TextCell "D7" '63.43' "D23"
Set-TextCell "E7" '  +1.94%  ' "E23"
Set-TextCell "E8" '  +0.09%  ' "E23"
Set-TextCell "E9" '  +0.66%  ' "E23"
Set-TextCell "D10" '0.0847' "D23"
Set-TextCell "E10" '  +0.36%  ' "E23"
Set-TextCell "E11" '  -0.03%  ' "E23"
Set-TextCell "D12" '15.96' "D23"
Set-TextCell "E12" '  +0.91%  ' "E23"
Set-TextCell "D13" '2.482.66' "D23"
Set-TextCell "E13" '  +2.70%  ' "E23"
Set-TextCell "D14" '21.96' "D23"
Set-TextCell "E14" '  -0.26%  ' "E23"
Set-TextCell "D15" '0.807' "D23"
Set-TextCell "E15" '  -0.17%  ' "E23"
Set-TextCell "D16" '5.48' "D23"
Set-TextCell "E16" '  -0.95%  ' "E23"
Set-TextCell "D17" '2.157.88' "D23"
Set-TextCell "E17" '  +3.27%  ' "E23"
Set-TextCell "D18" '39.562.01' "D23"
Set-TextCell "E18" '  +1.95%  ' "E23"
Set-TextCell "D19" '71.79' "D23"
Set-TextCell "E19" '  +0.26%  ' "E23"
Set-TextCell "D20" '6.10' "D23"
Set-TextCell "E20" '  -0.29%  ' "E23"
Set-TextCell "D21" '0.0₃0844' "D23"
Set-TextCell "E21" '  -0.31%  ' "E23"
Set-TextCell "D22" '227.75' "D23"
Set-TextCell "E22" '  -0.23%  ' "E23"
Set-TextCell "D24" '2.40' "D23"
Set-TextCell "E24" '  +3.04%  ' "E23"
Set-TextCell "E25" '  +1.17%  ' "E23"
Set-TextCell "B26" 'Monero' "B23"
Set-TextCell "C26" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' "C23"
Set-TextCell "D26" '172.75' "D23"
Set-TextCell "E26" '  +0.35%  ' "E23"
Set-TextCell "B27" 'Cosmos' "B23"
Set-TextCell "C27" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom' "C23"
Set-TextCell "D27" '9.63' "D23"
Set-TextCell "E27" '  -0.34%  ' "E23"
Set-TextCell "D28" '0.138' "D23"
Set-TextCell "E28" '  +0.56%  ' "E23"
Set-TextCell "D29" '19.71' "D23"
Set-TextCell "E29" '  +1.88%  ' "E23"
Set-TextCell "E30" '  -0.13%  ' "E23"
Set-TextCell "E31" '  +3.52%  ' "E23"
Set-TextCell "E32" '  +1.40%  ' "E23"
Set-TextCell "D33" '4.59' "D23"
Set-TextCell "E33" '  +0.27%  ' "E23"
Set-TextCell "D34" '4.67' "D23"
Set-TextCell "E34" '  -1.63%  ' "E23"
Set-TextCell "D35" '6.95' "D23"
Set-TextCell "E35" '  -3.33%  ' "E23"
Set-TextCell "D36" '0.0617' "D23"
Set-TextCell "E36" '  -0.16%  ' "E23"
Set-TextCell "D37" '2.40' "D23"
Set-TextCell "E37" '  +0.28%  ' "E23"
Set-TextCell "D38" '3.62' "D23"
Set-TextCell "E38" '  +2.78%  ' "E23"
Set-TextCell "B39" 'FTXToken' "B23"
Set-TextCell "C39" 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' "C23"
Set-TextCell "D39" '5.13' "D23"
Set-TextCell "E39" '  +23.67%  ' "E23"
Set-TextCell "B40" 'BinanceUSD' "B23"
Set-TextCell "C40" 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd' "C23"
Set-TextCell "D40" '1.00' "D23"
Set-TextCell "E40" '  +0.46%  ' "E23"
Set-TextCell "D41" '101.97' "D23"
Set-TextCell "E41" '  -0.27%  ' "E23"
Set-TextCell "D42" '0.0226' "D23"
Set-TextCell "E42" '  -0.48%  ' "E23"
Set-TextCell "B43" 'TrustWalletToken' "B23"
Set-TextCell "C43" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt' "C23"
Set-TextCell "D43" '1.24' "D23"
Set-TextCell "E43" '  +2.95%  ' "E23"
Set-TextCell "B44" 'InjectiveProtocol' "B23"
Set-TextCell "C44" 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj' "C23"
Set-TextCell "D44" '17.48' "D23"
Set-TextCell "E44" '  -3.51%  ' "E23"
Set-TextCell "D45" '1.512.07' "D23"
Set-TextCell "E45" '  -0.88%  ' "E23"
Set-TextCell "D46" '0.0922' "D23"
Set-TextCell "E46" '  +0.49%  ' "E23"
Set-TextCell "D48" '1.09' "D23"
Set-TextCell "E48" '  +1.25%  ' "E23"
Set-TextCell "D49" '7.76' "D23"
Set-TextCell "E49" '  +0.22%  ' "E23"
Set-TextCell "E50" '  +1.07%  ' "E23"
Set-TextCell "D51" '2.367.49' "D23"
Set-TextCell "E51" '  +2.74%  ' "E23"
